$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (E9): previously a plain boolean TRUE, now expressed as a TRUE() formula
$ws.Range("E9").Formula = "=TRUE()"

# New row 10 with a fresh dataset entry
$ws.Range("A10").Value = "11_14_21"
$ws.Range("B10").Value = 15
$ws.Range("D10").Value = 28.7
$ws.Range("E10").Value = $false
$ws.Range("E10").NumberFormat = $ws.Range("E9").NumberFormat
$ws.Range("F10").Value = 1000000

# Move the active selection as recorded in the saved view state
[void]$ws.Range("C14").Select()
